$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Header figures that changed value (same labels, new numbers)
# -----------------------------------------------------------------
$ws.Range("E11").Value = 401613      # Valor Mora total
$ws.Range("C13").Value = 4           # Cant. Trabajadores
$ws.Range("F13").Value = 6           # Cant. Periodos

# -----------------------------------------------------------------
# 2. Give row 23 the "closing" border style used by the last row of
#    the worker table (copy formatting only from the old last row,
#    row 35, before that row gets removed below).
# -----------------------------------------------------------------
$ws.Range("B35:J35").Copy() | Out-Null
$ws.Range("B23:J23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# -----------------------------------------------------------------
# 3. Replace the worker / period rows (16-23) with the new data set
# -----------------------------------------------------------------
$rows = @(
    @{ r=16; b="CC"; c="1047473239"; d="ALLISON YUSENIA MARINEZ ZURIQUE"; e="2402"; f=27733; g=1300000 },
    @{ r=17; b="CC"; c="1047473239"; d="ALLISON YUSENIA MARINEZ ZURIQUE"; e="2403"; f=52000; g=1300000 },
    @{ r=18; b="CC"; c="1143401687"; d="MARIA JOSE TROYA FONSECA";        e="2403"; f=52000; g=1423500 },
    @{ r=19; b="CC"; c="1047473239"; d="ALLISON YUSENIA MARINEZ ZURIQUE"; e="2404"; f=52000; g=1300000 },
    @{ r=20; b="CC"; c="1047473239"; d="ALLISON YUSENIA MARINEZ ZURIQUE"; e="2405"; f=52000; g=1300000 },
    @{ r=21; b="CC"; c="1047473239"; d="ALLISON YUSENIA MARINEZ ZURIQUE"; e="2406"; f=52000; g=1300000 },
    @{ r=22; b="CC"; c="52807924";   d="MARIA ALEXANDRA MOJICA OROZCO";   e="2508"; f=56940; g=1423500 },
    @{ r=23; b="CC"; c="1143370388"; d="LAURA PAOLA CASTRO CANTILLO";     e="2508"; f=56940; g=1423500 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("B$r").Value = $row.b
    $ws.Range("C$r").Value = $row.c
    $ws.Range("D$r").Value = $row.d
    $ws.Range("E$r").Value = $row.e
    $ws.Range("F$r").Value = $row.f
    $ws.Range("G$r").Value = $row.g
}

# -----------------------------------------------------------------
# 4. Remove the now-obsolete rows (old rows 24-35 of the previous
#    worker table); this also shifts the footer block up so it ends
#    up on rows 28-29 instead of 40-41.
# -----------------------------------------------------------------
$ws.Rows("24:35").Delete() | Out-Null

# -----------------------------------------------------------------
# 5. Column D ("Nombre Trabajador") is best-fit to its contents; the
#    shorter names now in use mean it should shrink to fit.
# -----------------------------------------------------------------
$ws.Columns("D").AutoFit() | Out-Null
